$d = $word.ActiveDocument

# --- Change 1: merge the three runs of the "Lydia - Provided data sources..."
# paragraph (week 1 section) into a single run. Doing a Find/Replace across the
# whole (already-contiguous) text forces the engine to re-emit it as one run.
$d.Content.Find.Execute(
    "Lydia - Provided data sources and descriptions, Questions to answer",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Lydia - Provided data sources and descriptions, Questions to answer", 2) | Out-Null

# --- Change 2: append a new "Work completed for Week 3:" section after the
# "Kelsey - Database" paragraph (still before the trailing blank paragraph).
$kelsey = $d.Paragraphs(14)
$kelsey.Range.InsertParagraphAfter()
$target = $d.Paragraphs(15)

$body =  '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>'
$body += '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Work completed for </w:t></w:r>'
$body += '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Week</w:t></w:r>'
$body += '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
$body += '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3</w:t></w:r>'
$body += '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r></w:p>'
$body += '<w:p><w:r><w:t>Ben – Updated ReadMe, Dashboard</w:t></w:r>'
$body += '<w:r><w:t xml:space="preserve"> - tableau</w:t></w:r></w:p>'
$body += '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Digamber</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$body += '<w:r><w:t xml:space="preserve"> – Updated ReadMe, Data Modeling</w:t></w:r></w:p>'
$body += '<w:p><w:r><w:t>Lydia – Updated ReadMe</w:t></w:r>'
$body += '<w:r><w:t xml:space="preserve">, </w:t></w:r>'
$body += '<w:r><w:t>Database, presentation outline</w:t></w:r></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml)
